# Update the "Förändrad" (Changed) date column (C) for rows 2-31 from
# 2026-02-21 (serial 46074) to 2026-02-22 (serial 46075).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
